# repull data, push all data, mean calculation
# Update the dSF column (F) with re-pulled values for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value  = -2
$ws.Range("F8").Value  = -5
$ws.Range("F9").Value  = -5
$ws.Range("F10").Value = -8
$ws.Range("F12").Value = -4
$ws.Range("F14").Value = 4
$ws.Range("F16").Value = -1
$ws.Range("F17").Value = -3
$ws.Range("F18").Value = -5
$ws.Range("F20").Value = 0
$ws.Range("F22").Value = 1
$ws.Range("F23").Value = 5
